# Renames the embedded logo images that appear in the document's
# header and footer:
#   - The BTec logo (currently "image1.jpg") becomes "image2.jpg"
#   - The Pearson logo (currently "image2.png") becomes "image1.png"
#
# Word.Sections(1).Headers/Footers index mapping in this document:
#   Headers.Item(1)  -> default header  (BTec logo,    currently image1.jpg)
#   Headers.Item(2)  -> first-page header (BTec logo,   currently image1.jpg)
#   Footers.Item(1)  -> default footer  (Pearson logo, currently image2.png)
#   Footers.Item(2)  -> first-page footer (Pearson logo, currently image2.png)

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers: BTec_Logo-Orange, image1.jpg -> image2.jpg ---
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
        $ishp = $hdr.Range.InlineShapes.Item($j)
        if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
            $ishp.Name = "image2.jpg"
        }
    }
}

# --- Footers: Pearson logo, image2.png -> image1.png ---
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
        $ishp = $ftr.Range.InlineShapes.Item($j)
        if ($ishp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $ishp.Name = "image1.png"
        }
    }
}
